# Cosmetic update to the "2018" ratings-history sheet: two new date
# columns (Jun_17, Jun_15) are inserted right after column A, pushing the
# existing Jun_13 column (B) and the rating-history column (C) two slots
# to the right (-> D, E respectively). The two freshly inserted columns
# are seeded with the same "UN" placeholder already used in the old B
# column, and their width matches the neighbouring column (8 characters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before the existing "Jun_13" column (B),
# shifting old B -> D and old C -> E while carrying their values/styles.
$ws.Columns("B:C").Insert()

# New header row entries for the freshly-opened columns (Jun_15 is
# written first so it lands before Jun_17 in the shared-string table,
# matching the order the dates were appended in).
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Placeholder data ("UN") for the new columns, mirroring column D (the
# old column B) for every data row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Match the new columns' width to the existing 8-character-wide column.
$ws.Columns("C:E").ColumnWidth = 7.17
